$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Distance of best" column header label from the three result tables,
# but keep the cell's existing (bold) formatting.
$ws.Range("E8").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("E32").ClearContents()

# Remove the stray "0.75 meters" note next to the first table.
$ws.Range("E16").ClearContents()

# Add the new motor, QBL5704-116-04-042, as a new row in the motor comparison
# table at the top of the sheet (Gear Ratio / dc supply voltage / armature
# inductance / no load speed / stall torque / Maximum current columns).
$ws.Range("A5").Value = "QBL5704-116-04-042"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = 36
$ws.Range("C5").Value = 0.001
$ws.Range("C5").NumberFormat = "0.00E+00"
$ws.Range("D5").Value = 5500
$ws.Range("E5").Value = 1.3
$ws.Range("F5").Value = 11

# Restore the view to the top of the sheet with F6 selected.
$ws.Range("F6").Select()
